$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'38.137.02"
$ws.Range("E2").Value = "  +3.13%  "
$ws.Range("D3").Value = "'2.061.72"
$ws.Range("E3").Value = "  +2.79%  "
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("D5").Value = "'230.50"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("E6").Value = "  +1.76%  "
$ws.Range("D7").Value = "'58.05"
$ws.Range("E7").Value = "  +6.55%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").Value = "'2.363.00"
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").Value = "'14.64"
$ws.Range("E13").Value = "  +3.62%  "
$ws.Range("D14").Value = "'20.68"
$ws.Range("E14").Value = "  +2.59%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'5.30"
$ws.Range("E15").Value = "  +3.84%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.754"
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("D17").Value = "'2.061.87"
$ws.Range("E17").Value = "  +3.11%  "
$ws.Range("D18").Value = "'38.068.40"
$ws.Range("E18").Value = "  +3.23%  "
$ws.Range("D19").Value = "'6.20"
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("D20").Value = "'69.70"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("D21").Value = "'0.0₃0831"
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("D22").Value = "'224.56"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").Value = "'2.25"
$ws.Range("E25").Value = "  +3.42%  "
$ws.Range("E26").Value = "  +2.18%  "
$ws.Range("D27").Value = "'165.91"
$ws.Range("D28").Value = "'0.134"
$ws.Range("E28").Value = "  +7.79%  "
$ws.Range("E29").Value = "  +2.18%  "
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("E32").Value = "  +1.52%  "
$ws.Range("E33").Value = "  +4.68%  "
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("E35").Value = "  +7.56%  "
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("D37").Value = "'6.06"
$ws.Range("E37").Value = "  +13.22%  "
$ws.Range("D38").Value = "'3.33"
$ws.Range("E38").Value = "  +6.39%  "
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").Value = "'98.53"
$ws.Range("E40").Value = "  +4.32%  "
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("D42").Value = "'1.484.22"
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "'0.0945"
$ws.Range("E43").Value = "  +2.95%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'16.86"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("E45").Value = "  +4.41%  "
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").Value = "'4.11"
$ws.Range("E47").Value = "  +18.97%  "
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("D50").Value = "'7.11"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").Value = "'2.250.65"
$ws.Range("E51").Value = "  +2.44%  "
